$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.214000701904297
$ws.Range("B1").Value = 2.587107419967651
$ws.Range("C1").Value = 4.282364845275879
$ws.Range("D1").Value = 2.042232275009155
$ws.Range("E1").Value = 1.16678786277771
